$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from 2023-10-22 (45221) to 2023-10-25 (45224)
foreach ($row in 2..11) {
    $ws.Cells.Item($row, 3).Value = 45224
}
